$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 4; this shifts the former row 4 (Petro Matad) down to row 5
$ws.Rows.Item(4).Insert()

# --- Row 2: update company index + metrics ---
$ws.Cells.Item(2,1).Value = "Isle of Man"
$ws.Cells.Item(2,2).NumberFormat = "@"
$ws.Cells.Item(2,2).Value = "3"
$ws.Cells.Item(2,2).Style = "Normal"
$ws.Cells.Item(2,3).Value = "Oil/Gas (Production and Exploration)"
$ws.Cells.Item(2,7).Value = -1043.181818181818
$ws.Cells.Item(2,8).Value = -1043.181818181818
$ws.Cells.Item(2,9).Value = -961.0337744215246
$ws.Cells.Item(2,10).Value = -961.0337744215246
$ws.Cells.Item(2,11).Value = -20.943
$ws.Cells.Item(2,12).Value = -951.9545454545456
$ws.Cells.Item(2,13).Value = 0
$ws.Cells.Item(2,14).Value = 0
$ws.Cells.Item(2,15).Value = -0
$ws.Cells.Item(2,16).Value = 0
$ws.Cells.Item(2,17).Value = 0
$ws.Cells.Item(2,18).Value = -0
$ws.Cells.Item(2,19).Value = 0
$ws.Cells.Item(2,21).Value = 15.105
$ws.Cells.Item(2,22).Value = 0.08832816794339513
$ws.Cells.Item(2,23).Value = -0.1024952015355086
$ws.Cells.Item(2,24).Value = 0.05873265225172074
$ws.Cells.Item(2,25).Value = -0.1612278537872294
$ws.Cells.Item(2,26).Value = 0.0002473669052640273
$ws.Cells.Item(2,27).Value = -0.1075159546359904
$ws.Cells.Item(2,28).Value = 0.05872481633742673
$ws.Cells.Item(2,29).Value = -0.1662407709734171
$ws.Cells.Item(2,30).Value = 0.089
$ws.Cells.Item(2,31).Value = 0.04371518636770249
$ws.Cells.Item(2,32).Value = 0.1327151863677025
$ws.Cells.Item(2,33).Value = -14.9722848136323
$ws.Cells.Item(2,34).Value = 0.0007754650042987857
$ws.Cells.Item(2,35).Value = 0.001400943550563413
$ws.Cells.Item(2,36).Value = -0.0959529867234326
$ws.Cells.Item(2,37).Value = -0.1880285623992833
$ws.Cells.Item(2,38).Value = 0.097
$ws.Cells.Item(2,39).Value = -0.247
$ws.Cells.Item(2,40).Value = -0.004258780744568858
$ws.Cells.Item(2,41).Value = -218.3711340206185
$ws.Cells.Item(2,42).Value = 0.7164458232190782
$ws.Cells.Item(2,43).Value = 85.7570850202429

# --- Row 3 (TomCo Energy): update metrics ---
$ws.Cells.Item(3,1).Value = "Isle of Man"
$ws.Cells.Item(3,2).Value = "TomCo Energy Plc (AIM:TOM)"
$ws.Cells.Item(3,3).Value = "Oil/Gas (Production and Exploration)"
$ws.Cells.Item(3,11).Value = -0.803
$ws.Cells.Item(3,13).Value = -0
$ws.Cells.Item(3,14).Value = -0
$ws.Cells.Item(3,15).Value = 0
$ws.Cells.Item(3,16).Value = -0
$ws.Cells.Item(3,17).Value = -0
$ws.Cells.Item(3,18).Value = 0
$ws.Cells.Item(3,19).Value = 0
$ws.Cells.Item(3,21).Value = 0.9350000000000001
$ws.Cells.Item(3,22).Value = 0.1004296455424275
$ws.Cells.Item(3,23).Value = -0.06636363636363637
$ws.Cells.Item(3,24).Value = 0.05893354543001046
$ws.Cells.Item(3,25).Value = -0.1252971817936469
$ws.Cells.Item(3,26).Value = 0
$ws.Cells.Item(3,27).Value = -0.06365844973089513
$ws.Cells.Item(3,28).Value = 0.05911620255356766
$ws.Cells.Item(3,29).Value = -0.1227746522844628
$ws.Cells.Item(3,30).Value = 0
$ws.Cells.Item(3,31).Value = 0.04371518636770249
$ws.Cells.Item(3,32).Value = 0.04371518636770249
$ws.Cells.Item(3,33).Value = -0.8912848136322976
$ws.Cells.Item(3,34).Value = 0.004673563979306736
$ws.Cells.Item(3,35).Value = 0.003430332970283722
$ws.Cells.Item(3,36).Value = -0.1058694579756709
$ws.Cells.Item(3,37).Value = -0.07547686598972433
$ws.Cells.Item(3,38).Value = 0
$ws.Cells.Item(3,39).Value = -0.002
$ws.Cells.Item(3,40).Value = -0
$ws.Cells.Item(3,42).Value = 1.207702999501758
$ws.Cells.Item(3,43).Value = 396

# Row 3 no longer has an ebit_interest_expenses (AO) value; remove it
$ws.Cells.Item(3,41).ClearContents()

# --- Row 4 (new): Bahamas Petroleum Company plc (AIM:BPC) ---
$ws.Cells.Item(4,1).Value = "Isle of Man"
$ws.Cells.Item(4,2).Value = "Bahamas Petroleum Company plc (AIM:BPC)"
$ws.Cells.Item(4,3).Value = "Oil/Gas (Production and Exploration)"
$ws.Cells.Item(4,11).Value = -5.34
$ws.Cells.Item(4,13).Value = -0
$ws.Cells.Item(4,14).Value = -0
$ws.Cells.Item(4,15).Value = 0
$ws.Cells.Item(4,16).Value = -0
$ws.Cells.Item(4,17).Value = -0
$ws.Cells.Item(4,18).Value = 0
$ws.Cells.Item(4,19).Value = 0
$ws.Cells.Item(4,21).Value = 12.1
$ws.Cells.Item(4,22).Value = 0.08551236749116607
$ws.Cells.Item(4,23).Value = -0.1024952015355086
$ws.Cells.Item(4,24).Value = 0.05873265225172074
$ws.Cells.Item(4,25).Value = -0.1612278537872294
$ws.Cells.Item(4,26).Value = 0
$ws.Cells.Item(4,27).Value = -0.1075159546359904
$ws.Cells.Item(4,28).Value = 0.05872481633742673
$ws.Cells.Item(4,29).Value = -0.1662407709734171
$ws.Cells.Item(4,30).Value = 0.089
$ws.Cells.Item(4,31).Value = 0
$ws.Cells.Item(4,32).Value = 0.089
$ws.Cells.Item(4,33).Value = -12.011
$ws.Cells.Item(4,34).Value = 0.0006285799038060866
$ws.Cells.Item(4,35).Value = 0.00138437368756708
$ws.Cells.Item(4,36).Value = -0.09275691371467845
$ws.Cells.Item(4,37).Value = -0.2301442832780854
$ws.Cells.Item(4,38).Value = 0.097
$ws.Cells.Item(4,39).Value = 0.027
$ws.Cells.Item(4,40).Value = -0.01692015209125475
$ws.Cells.Item(4,41).Value = -54.5360824742268
$ws.Cells.Item(4,42).Value = 2.283460076045627
$ws.Cells.Item(4,43).Value = -195.925925925926

# --- Row 5 (was row 4): Petro Matad Limited - refreshed metrics ---
$ws.Cells.Item(5,1).Value = "Isle of Man"
$ws.Cells.Item(5,2).Value = "Petro Matad Limited (DB:HA3)"
$ws.Cells.Item(5,3).Value = "Oil/Gas (Production and Exploration)"
$ws.Cells.Item(5,7).Value = -790.9090909090909
$ws.Cells.Item(5,8).Value = -790.9090909090909
$ws.Cells.Item(5,9).Value = -686.3636363636364
$ws.Cells.Item(5,10).Value = -686.3636363636364
$ws.Cells.Item(5,11).Value = -14.8
$ws.Cells.Item(5,12).Value = -672.7272727272727
$ws.Cells.Item(5,13).Value = -0
$ws.Cells.Item(5,14).Value = -0
$ws.Cells.Item(5,15).Value = 0
$ws.Cells.Item(5,16).Value = -0
$ws.Cells.Item(5,17).Value = -0
$ws.Cells.Item(5,18).Value = 0
$ws.Cells.Item(5,19).Value = 0
$ws.Cells.Item(5,21).Value = 2.07
$ws.Cells.Item(5,22).Value = 0.1024752475247525
$ws.Cells.Item(5,23).Value = -0.4668769716088328
$ws.Cells.Item(5,24).Value = 0.05870157987992347
$ws.Cells.Item(5,25).Value = -0.5255785514887563
$ws.Cells.Item(5,26).Value = 0.0007882479398065209
$ws.Cells.Item(5,27).Value = -0.5410247223217485
$ws.Cells.Item(5,28).Value = 0.05870157987992347
$ws.Cells.Item(5,29).Value = -0.5997263022016719
$ws.Cells.Item(5,30).Value = 0
$ws.Cells.Item(5,31).Value = 0
$ws.Cells.Item(5,32).Value = 0
$ws.Cells.Item(5,33).Value = -2.07
$ws.Cells.Item(5,34).Value = 0
$ws.Cells.Item(5,35).Value = 0
$ws.Cells.Item(5,36).Value = -0.1141753998896856
$ws.Cells.Item(5,37).Value = -0.1324376199616123
$ws.Cells.Item(5,38).Value = 0
$ws.Cells.Item(5,39).Value = -0.272
$ws.Cells.Item(5,40).Value = -0
$ws.Cells.Item(5,42).Value = 0.1389261744966443
$ws.Cells.Item(5,43).Value = 55.51470588235293
